# Boid simulation workbook: add a test for moving a boid one step.
#
# The sheet already has a "before step" adjustments block (rows 27-32) and a
# "after step" / bounds block (rows 36-41, with the totals in row 41). This
# change:
#   1. Inserts a new "self-velocity" row before the existing totals row in the
#      rows 36-41 block (pushing the totals formula down to row 42), and
#      extends the totals SUM() ranges to include it.
#   2. Adds a new "New Boid " row below the totals that adds the bounds
#      constant (-6 / 6) back to the new totals, giving the moved boid's
#      final position.
#   3. Goes back and adds the equivalent missing "self-velocity" row to the
#      earlier adjustments block (rows 27-32), extending those SUM() ranges
#      too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "after step" block: insert the self-velocity row above the totals ---

# Row 41 currently holds the totals formula; push it down to row 42 so we can
# use row 41 for the new self-velocity data.
$ws.Rows("41:41").Insert()

$ws.Range("A41").Value = "self-velocity"
$ws.Range("C41").Value = -1
$ws.Range("D41").Value = 2

# Extend the (now shifted) totals formulas to include the new row.
$ws.Range("C42").Formula = "=SUM(C37:C41)"
$ws.Range("D42").Formula = "=SUM(D37:D41)"

# New row: add the bounds constant back to get the moved boid's new position.
$ws.Range("A45").Value = "New Boid "
$ws.Range("C45").Formula = "=-6+C42"
$ws.Range("D45").Formula = "=6+D42"

# --- earlier adjustments block: add the matching self-velocity row ---

$ws.Range("A31").Value = "self-velocity"
$ws.Range("C31").Value = 3
$ws.Range("D31").Value = 5

$ws.Range("C32").Formula = "=SUM(C28:C31)"
$ws.Range("D32").Formula = "=SUM(D28:D31)"

# The new self-velocity row pushed the chart's anchor down by one row
# (13pt); grow the chart to follow it, since it's set to move/size with
# cells.
$chart = $ws.ChartObjects().Item(1)
$chart.Height = $chart.Height + 13

# Update the current selection, as left by the editing session.
$ws.Range("D33").Select()
